# Updates cryptos list data (price + 1h volume change) to refreshed values,
# and corrects the ranking order for two coin pairs (LEO/Toncoin, Hedera/Cosmos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.218.14"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "4.028.09"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'528.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'150.91"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.701"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +12.05%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.749"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D11").Value = "'0.0000327"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'49.79"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.83%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'10.75"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "4.670.05"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "4.092.29"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'14.10"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'20.65"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.77%  "
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "72.133.79"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'433.97"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'98.04"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'3.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'4.19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.95%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'14.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'11.26"
$ws.Range("D26").ClearFormats()
$ws.Range("D27").Value = "'10.74"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.06%  "
$ws.Range("E27").ClearFormats()
$ws.Range("B28").Value = "Toncoin"
$ws.Range("B28").ClearFormats()
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C28").ClearFormats()
$ws.Range("D28").Value = "'3.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +19.85%  "
$ws.Range("E28").ClearFormats()
$ws.Range("B29").Value = "LEO"
$ws.Range("B29").ClearFormats()
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C29").ClearFormats()
$ws.Range("D29").Value = "'5.86"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.39%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'36.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'7.48"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.86%  "
$ws.Range("E31").ClearFormats()
$ws.Range("B32").Value = "Cosmos"
$ws.Range("B32").ClearFormats()
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C32").ClearFormats()
$ws.Range("D32").Value = "'13.47"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E32").ClearFormats()
$ws.Range("B33").Value = "Hedera"
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C33").ClearFormats()
$ws.Range("D33").Value = "'0.131"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'681.03"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'48.13"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +17.71%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'65.79"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.448"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "0.0₃0827"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -9.56%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "  -7.68%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'3.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.69%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "  +0.88%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'10.33"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.08%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "  -4.04%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'3.39"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'3.02"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.57%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.000269"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'3.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.38%  "
$ws.Range("E51").ClearFormats()
